$wb = $excel.ActiveWorkbook

# --- Sheet "BasePath_Directory" (sheet1) -------------------------------
$ws1 = $wb.Worksheets.Item("BasePath_Directory")

# Move existing rows 12/13/14 down to 16/19/20 to make room for the new
# "SystemPopup" entry, then insert the new row 15 data.
$ws1.Range("A16").Value() = $ws1.Range("A12").Value()
$ws1.Range("B16").Value() = $ws1.Range("B12").Value()

$ws1.Range("A19").Value() = $ws1.Range("A13").Value()
$ws1.Range("B19").Value() = $ws1.Range("B13").Value()

$ws1.Range("A20").Value() = $ws1.Range("A14").Value()
$ws1.Range("B20").Value() = $ws1.Range("B14").Value()

$ws1.Range("A12:B12").ClearContents()
$ws1.Range("A13:B13").ClearContents()
$ws1.Range("A14:B14").ClearContents()

# --- Sheet "BasePath_BP_File" (sheet2) ----------------------------------
$ws2 = $wb.Worksheets.Item("BasePath_BP_File")

# New row 16: Id=1000, Directory_Table_Id=999, BP_File_Name="SystemPopup"
# (entered first so it claims the earlier shared-string slot)
$ws2.Range("A16").Value() = 1000
$ws2.Range("B16").Value() = 999
$ws2.Range("C16").Value() = "SystemPopup"

# New row 15 on sheet1: Id=999, Directory="Widget/BuiltInWidget/SystemWidget"
$ws1.Range("A15").Value() = 999
$ws1.Range("B15").Value() = "Widget/BuiltInWidget/SystemWidget"

# Page setup (A4, portrait) for this sheet
$ws1.PageSetup.PaperSize = 9
$ws1.PageSetup.Orientation = 1

# --- Selection / active sheet ------------------------------------------
# Final active sheet is BasePath_Directory, selection B15.
# Sheet2's selection becomes C20 (no longer the active tab).
$null = $ws2.Range("C20").Select()
$null = $ws1.Range("B15").Select()
